# Generate Report for Handback
#
# Row 7 in both the "zh-cn" and "de-de" sheets corresponds to the
# a5c5a14d-4454-4df0-a5ab-5ee4028c9d3f handback doc. A new handback was
# received for that doc (for both target languages), so the "Latest
# Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns get populated and an error is raised because the handback
# doesn't match the very latest source version.

$wb = $excel.ActiveWorkbook

$docId   = "a5c5a14d-4454-4df0-a5ab-5ee4028c9d3f"
$docLink = "$docId.md"

$errorMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/006bb043a47460cd6cac2ca97ff5427f824d0a40/e2e/$docId.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/550f0d84ec416eac11db159d538e6bd10e82b24b/e2e/$docId.md."

# ---------------------------------------------------------------------
# zh-cn sheet, row 7
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# I: Latest Target File -> new hyperlink to the handed-back doc in the
# zh-cn target repo.
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/550f0d84ec416eac11db159d538e6bd10e82b24b/e2e/$docLink",
    "",
    "",
    $docLink
) | Out-Null

# J: Latest Handback File
$wsZhCn.Range("J7").Value = "$docId.b86efd542a39abff032e68b2350cefab549de52e.zh-cn.xlf"

# K: Latest Handback DateTime
$wsZhCn.Range("K7").Value = "2016-08-23 10:56:49"

# P: Error Detail
$wsZhCn.Range("P7").Value = $errorMsg

# ---------------------------------------------------------------------
# de-de sheet, row 7
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# I: Latest Target File -> new hyperlink to the handed-back doc in the
# de-de target repo.
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/550f0d84ec416eac11db159d538e6bd10e82b24b/e2e/$docLink",
    "",
    "",
    $docLink
) | Out-Null

# J: Latest Handback File
$wsDeDe.Range("J7").Value = "$docId.b86efd542a39abff032e68b2350cefab549de52e.de-de.xlf"

# K: Latest Handback DateTime
$wsDeDe.Range("K7").Value = "2016-08-23 10:56:56"

# P: Error Detail
$wsDeDe.Range("P7").Value = $errorMsg
